$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 3).Value = "'25"
$ws.Cells.Item($row, 3).Style = "Normal"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
